## issue #5: stock data output to json file
##
## The "股票" (stock) sheet gains a new "property_category" column (constant
## value "stock") inserted right after the "total" column and before the
## "date" column. That pushes date / legislator_name / legislator_id one
## column to the right (H->I, I->J, J->K).
##
## Also: a handful of numeric-looking values in that sheet were stored as
## text with (half-width/full-width) thousands separators, e.g. "2,275,950"
## or "76，330". Those get normalised to plain digit strings ("2275950",
## "76330", ...) while remaining text cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# --- Insert the new "property_category" column at H, shifting the old
#     H:J (date / legislator_name / legislator_id) columns to I:K. The
#     freshly inserted column inherits the formatting of the column that
#     used to be there (header style on row 1, data style elsewhere).
$ws.Columns("H:H").Insert()

$ws.Range("H1").Value = "property_category"
$ws.Range("H2:H39").Value = "stock"

# --- Normalise the comma-formatted numeric text values so they no longer
#     contain thousands separators, while keeping them stored as text
#     (not auto-converted to numbers) by briefly marking the cells as
#     Text-formatted, assigning the value, then restoring the style so the
#     look of the cell (General format) is unchanged.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D21") "5599"
Set-TextValue $ws.Range("D29") "26773"
Set-TextValue $ws.Range("G7") "2275950"
Set-TextValue $ws.Range("G11") "1000000"
Set-TextValue $ws.Range("G12") "76330"
Set-TextValue $ws.Range("G15") "2740"
Set-TextValue $ws.Range("G21") "55990"
Set-TextValue $ws.Range("G27") "3090"
Set-TextValue $ws.Range("G28") "4050"
Set-TextValue $ws.Range("G31") "203960"
Set-TextValue $ws.Range("G33") "236640"
Set-TextValue $ws.Range("G38") "3450"
